$d = $word.ActiveDocument
$d.Content.Find.Execute(" Bash, Beautiful Soup,", $true, $false, $false, $false, $false,
                         $true, 1, $false, " Bash,", 2)
